$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 21 with the 2025-09-05 portfolio data.
# The date column holds plain text like "2025-09-05" (matching the existing
# rows), so force text formatting first to stop Excel from auto-converting
# it into a date serial number, then drop the formatting again so the new
# cell ends up unstyled just like its neighbours.
$dateCell = $ws.Range("A21")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-09-05"
$dateCell.ClearFormats()

$ws.Range("B21").Value = 57.86999893188477
$ws.Range("C21").Value = 691.7000122070312
$ws.Range("D21").Value = 329.1499938964844
